# Rename the existing "Datos" sheet to "caja", then duplicate it (data,
# formatting, column widths, etc.) into a new sheet named "post" placed
# right after "caja". Finally restore "caja" as the active sheet so the
# workbook's activeTab stays at 0, matching the original view state.

$wb = $excel.ActiveWorkbook

$caja = $wb.Worksheets.Item(1)
$caja.Name = "caja"

# Worksheet.Copy duplicates the sheet (cells, styles, column widths, ...)
# and inserts the copy right after the $caja sheet.
$caja.Copy([System.Reflection.Missing]::Value, $caja)

$post = $wb.Worksheets.Item(2)
$post.Name = "post"

# Keep the first sheet ("caja") selected/active, as in the source workbook.
$caja.Activate()
